$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 326.8
Write-Output $ws.Range("H9").Value
